# Daily attendance processing - 2025-11-30 20:28:08
# Normalizes the ordering of names/emails in the "Recorded By" column (G)
# so that "System" is listed before the associated email address.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "system, backup@backdoor.com, System" = "system, System, backup@backdoor.com"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
